# Fruta / hortaliza, semanal
# This weekly refresh re-sorts/re-pulls the daily rows (2-35) for this
# subset ("Hortaliza, Vega Modelo de Temuco - Achicoria"). The columns that
# identify the market/category (A,B,C,E,F,G,H,R) stay the same for every
# row, but the per-record columns (Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad de comercializacion,
# Origen, Precio $/Kg, Kg o Unidades) get reassigned row by row following
# the mapping below (destination row -> source row, both 1-based sheet rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Destination row -> source row mapping (rows 2..35 of the data table).
$map = @{2=6; 3=31; 4=23; 5=26; 6=12; 7=11; 8=5; 9=33; 10=7; 11=22; 12=10; 13=34; 14=24; 15=8; 16=30; 17=18; 18=14; 19=21; 20=20; 21=35; 22=13; 23=15; 24=16; 25=32; 26=9; 27=3; 28=19; 29=2; 30=27; 31=28; 32=29; 33=4; 34=17; 35=25}

# Columns that vary per record (others stay constant across all rows).
$cols = @(4, 9, 10, 11, 12, 13, 14, 15, 16, 17)   # D, I, J, K, L, M, N, O, P, Q

$firstRow = 2
$lastRow = 35

# 1) Snapshot the original values for all the columns we are about to move,
#    so that writes to destination rows do not clobber values we still
#    need to read as a source for a later row.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write back each destination row using the snapshot of its source row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $map[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $srcVals[$c]
    }
}
